$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.306.78'
$ws.Range("E2").Value = '  -6.44%  '
$ws.Range("D3").Value = '3.275.87'
$ws.Range("E3").Value = '  -7.72%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '180.85'
$ws.Range("E5").Value = '  -11.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '517.47'
$ws.Range("E6").Value = '  -7.08%  '
$ws.Range("E7").Value = '  -1.39%  '
$ws.Range("D8").Value = '3.270.53'
$ws.Range("E8").Value = '  -7.73%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.616'
$ws.Range("E10").Value = '  -7.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.24'
$ws.Range("E11").Value = '  -8.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.130'
$ws.Range("E12").Value = '  -9.77%  '
$ws.Range("E13").Value = '  -7.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.07'
$ws.Range("E14").Value = '  -9.15%  '
$ws.Range("D15").Value = '3.802.76'
$ws.Range("E15").Value = '  -7.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.117'
$ws.Range("E16").Value = '  -5.94%  '
$ws.Range("D17").Value = '3.280.29'
$ws.Range("E17").Value = '  -7.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.64'
$ws.Range("E18").Value = '  -6.43%  '
$ws.Range("D19").Value = '63.158.75'
$ws.Range("E19").Value = '  -6.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.91'
$ws.Range("E20").Value = '  -9.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.945'
$ws.Range("E21").Value = '  -9.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '369.96'
$ws.Range("E22").Value = '  -6.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.31'
$ws.Range("E23").Value = '  -8.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.69'
$ws.Range("E24").Value = '  -9.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.80'
$ws.Range("E25").Value = '  -4.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.81'
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.98'
$ws.Range("E27").Value = '  -1.83%  '
$ws.Range("E28").Value = '  -7.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.33'
$ws.Range("E29").Value = '  -7.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.29'
$ws.Range("E30").Value = '  -8.08%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.39'
$ws.Range("E31").Value = '  -8.51%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '633.29'
$ws.Range("E32").Value = '  -9.88%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.66'
$ws.Range("E33").Value = '  -10.23%  '
$ws.Range("E34").Value = '  -6.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.106'
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.69'
$ws.Range("E36").Value = '  -7.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.391'
$ws.Range("E38").Value = '  -5.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.00'
$ws.Range("E39").Value = '  -12.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = '2.984.60'
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.124'
$ws.Range("E42").Value = '  -5.04%  '
$ws.Range("D43").Value = '0.0₃0650'
$ws.Range("E43").Value = '  -9.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.66'
$ws.Range("E45").Value = '  -15.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.58'
$ws.Range("E46").Value = '  -5.07%  '
$ws.Range("E47").Value = '  -4.34%  '
$ws.Range("E48").Value = '  +5.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.125'
$ws.Range("E49").Value = '  -3.35%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.89'
$ws.Range("E50").Value = '  -4.31%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.44'
$ws.Range("E51").Value = '  -21.26%  '
